# Update symbol list (prices and volume percentages) per the Feb 11 2023 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (kept as text, matching the sheet's
# existing convention of storing Price/Volume columns as strings).
$updates = [ordered]@{
    "D2" = "308.29"
    "E2" = "0.52%"
    "D3" = "40.75"
    "E3" = "2.44%"
    "D4" = "5.121"
    "E4" = "0.07%"
    "D5" = "0.07610"
    "E5" = "-1.25%"
    "E6" = "-0.20%"
    "E7" = "0.45%"
    "D8" = "0.9008"
    "E8" = "2.37%"
    "D9" = "0.1100"
    "E9" = "9.58%"
    "D10" = "0.1768"
    "E10" = "1.24%"
    "D11" = "0.09134"
    "E11" = "2.49%"
    "D12" = "0.04167"
    "E12" = "-4.98%"
    "D13" = "0.1050"
    "E13" = "-0.49%"
    "D14" = "0.001254"
    "E14" = "-0.05%"
    "D15" = "0.005801"
    "E15" = "-0.68%"
    "E16" = "0.02%"
    "D17" = "4.256"
    "E17" = "0.32%"
    "E18" = "-0.93%"
    "D19" = "6.576"
    "E19" = "-6.17%"
    "D20" = "0.1365"
    "E20" = "1.92%"
    "D21" = "0.2681"
    "E21" = "-10.70%"
    "D22" = "0.04064"
    "E22" = "-2.24%"
    "E23" = "1.83%"
    "D24" = "0.004088"
    "E24" = "0.53%"
    "D25" = "0.0001300"
    "E25" = "6.53%"
    "D38" = "0.02374"
    "E38" = "1.83%"
    "D39" = "0.05183"
    "E39" = "0.76%"
    "D40" = "0.007754"
    "E40" = "-2.14%"
    "D41" = "0.1300"
    "E41" = "-1.58%"
    "D42" = "0.006755"
    "E42" = "6.74%"
    "D43" = "0.001951"
    "E43" = "0.01%"
    "D44" = "0.007941"
    "E44" = "-6.52%"
    "D45" = "0.3341"
    "E45" = "9.30%"
    "D46" = "0.00007017"
    "E46" = "7.76%"
    "D47" = "0.00000000750"
    "E47" = "-0.02%"
    "D48" = "0.02962"
    "E48" = "538.83%"
    "D49" = "0.004202"
    "E49" = "-40.02%"
    "D50" = "0.00002101"
    "E50" = "-0.02%"
    "E51" = "-0.02%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Prefix with an apostrophe so Excel stores the numeric/percent-looking
    # text as a literal string instead of converting it to a number.
    $cell.Value = "'" + $updates[$addr]
    # Reset to the default style so no new number-format style gets attached
    # to the cell (keeps formatting identical to the original file).
    $cell.Style = "Normal"
}
